# Apply "Add 2022-05-20 data" update to Fonds de solidarite volet 1
# regional / classe effectif dataset. Updates columns C (nombre_aides),
# D (nombre_entreprises) and E (montant_total) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> C, D, E values (nombre_aides, nombre_entreprises, montant_total)
$updates = @(
    @{ Row = 91;  C = 151114;  D = 24834;  E = 482250128 },
    @{ Row = 92;  C = 409043;  D = 70904;  E = 1593918662 },
    @{ Row = 93;  C = 209529;  D = 34261;  E = 1308143776 },
    @{ Row = 94;  C = 94174;   D = 13795;  E = 916746241 },
    @{ Row = 95;  C = 50737;   D = 6982;   E = 931249783 },
    @{ Row = 96;  C = 17250;   D = 2565;   E = 789924340 },
    @{ Row = 104; C = 135231;  D = 23286;  E = 272154522 },
    @{ Row = 114; C = 3800;    D = 699;    E = 9102747 },
    @{ Row = 118; C = 977;     D = 161;    E = 11823184 },
    @{ Row = 131; C = 75582;   D = 15099;  E = 307221165 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 3).Value = $u.C
    $ws.Cells.Item($r, 4).Value = $u.D
    $ws.Cells.Item($r, 5).Value = $u.E
}
